# Applies the update described by the commit "Atualizado por script em 11-11-2023 14:45"
# to the Persian Gulf Pro League 2023-2024 betting-odds worksheet:
#   1) Rows 35 and 36 (match records) have their F:V content swapped.
#   2) Rows 50 and 51 (match records) have their F:V content swapped.
#   3) Two new match rows (69 and 70) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowValues($row) {
    $vals = @()
    foreach ($c in $cols) {
        $vals += , ($ws.Range($c + $row).Value())
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowValues $rowA
    $b = Get-RowValues $rowB
    Set-RowValues $rowA $b
    Set-RowValues $rowB $a
}

# 1) Swap the match details between rows 35 and 36
Swap-Rows 35 36

# 2) Swap the match details between rows 50 and 51
Swap-Rows 50 51

# 3) Append the two new rows at the bottom (69 and 70), copying formatting
#    from the last existing data row (68) so styles (borders/number formats)
#    stay consistent with the rest of the table.
$ws.Range("A68:V68").Copy()
$ws.Range("A69:V70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row69 = @(68,"iran","persian-gulf-pro-league","2023-2024",45241.54166666666,"Sanat Naft",0,"Havadar SC",2,2.22,"10/11/2023 01:13",2.51,"11/11/2023 12:08",2.75,"10/11/2023 01:13",2.75,"11/11/2023 12:08",3.32,"10/11/2023 01:13",3.23,"11/11/2023 12:08","https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sanat-naft-havadar-sc/rq2gVqNq/")

$row70 = @(69,"iran","persian-gulf-pro-league","2023-2024",45241.57638888889,"Esteghlal F.C.",2,"Tractor",0,1.94,"10/11/2023 02:13",1.99,"11/11/2023 13:46",2.96,"10/11/2023 02:13",2.81,"11/11/2023 13:47",3.81,"10/11/2023 02:13",4.72,"11/11/2023 13:47","https://www.betexplorer.com/football/iran/persian-gulf-pro-league/esteghlal-teh-tractor/tzGRZPES/")

for ($i = 0; $i -lt $allCols.Length; $i++) {
    $ws.Range($allCols[$i] + "69").Value = $row69[$i]
    $ws.Range($allCols[$i] + "70").Value = $row70[$i]
}
